# "add last guess and game end"
#
# The "Turn End" packet block (row 38, "guest result") gets a typo fix
# ("guest result" -> "guess result") and a new line documenting the
# "guess word or char" / "p or python" field right below it. Everything
# that used to start at row 40 ("Disqualify" section onward) shifts down
# by one row to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 39, pushing "Disqualify" (old row 40) and
# everything after it down by one.
$ws.Rows("39").Insert()

# Fix the typo on the existing row (old "guest result" -> "guess result").
$ws.Range("D38").Value = "guess result"

# Populate the newly inserted row 39 with the new field description.
$ws.Range("B39").Value = "String"
$ws.Range("D39").Value = "guess word or char"
$ws.Range("C39").Value = "p or python"

# Reposition the view roughly where the author left it.
$ws.Range("C40").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
